# Daily update at 8 AM UTC
# Adds the new day row (row 74) to the Wins Over Time sheet and moves the
# "last row" date formatting (no time component) from the old last row
# (73) onto the new last row (74).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 73 (the old last row) goes back to the regular date/time format
# used by every other data row.
$ws.Range("A73").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New day of data.
$ws.Range("A74").Value = 45814
$ws.Range("B74").Value = 316
$ws.Range("C74").Value = 313
$ws.Range("D74").Value = 318

# Row 74 becomes the new "last row" with the date-only format.
$ws.Range("A74").NumberFormat = "YYYY-MM-DD"
